# Fruta / hortaliza, semanal
# Two new weekly price records are inserted into the dataset (Ciboulette,
# Femacal de La Calera), one before the existing row 85 and one before the
# existing row 258 (1-based sheet rows). Inserting a row shifts every row
# below it down by one, which reproduces the "every D/J (and occasionally
# K/L/M/P) value equals the previous row's old value" pattern seen in the
# diff, all the way down to the bottom of the table, where the former last
# row (44160 / 230) ends up unchanged two rows further down (old R271 -> R273).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert #1: new record before current row 85 -------------------------
$ws.Rows("85:85").Insert()

$ws.Range("A85").Value = 3
$ws.Range("B85").Value = "Femacal de La Calera"
$ws.Range("C85").Value = "Coquimbo"
$ws.Range("D85").Value2 = 44614
$ws.Range("E85").Value = 5
$ws.Range("F85").Value = 100112039
$ws.Range("G85").Value = "Ciboulette"
$ws.Range("H85").Value = "Sin especificar"
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 230
$ws.Range("K85").Value = 1500
$ws.Range("L85").Value = 1500
$ws.Range("M85").Value = 1500
$ws.Range("N85").Value = "$/docena de atados"
$ws.Range("O85").Value = "Provincia de Quillota"
$ws.Range("P85").Value = 500
$ws.Range("Q85").Value = 3
$ws.Range("R85").Value = "Hortaliza"

# --- Insert #2: new record before current row 258 ------------------------
# (row 258 before this second insert is the old row 257, since everything
# from row 85 down already shifted by one after insert #1)
$ws.Rows("258:258").Insert()

$ws.Range("A258").Value = 3
$ws.Range("B258").Value = "Femacal de La Calera"
$ws.Range("C258").Value = "Coquimbo"
$ws.Range("D258").Value2 = 44615
$ws.Range("E258").Value = 5
$ws.Range("F258").Value = 100112039
$ws.Range("G258").Value = "Ciboulette"
$ws.Range("H258").Value = "Sin especificar"
$ws.Range("I258").Value = "Primera"
$ws.Range("J258").Value = 160
$ws.Range("K258").Value = 1500
$ws.Range("L258").Value = 1500
$ws.Range("M258").Value = 1500
$ws.Range("N258").Value = "$/docena de atados"
$ws.Range("O258").Value = "Provincia de Quillota"
$ws.Range("P258").Value = 500
$ws.Range("Q258").Value = 3
$ws.Range("R258").Value = "Hortaliza"
